$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (week 37): Friday (H26) is now marked as "done" for the week, so it
# gets the same fill/border formatting as the rest of that row (D26:G26),
# and the day count in I26 goes from 4 to 5.
$ws.Range("G26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I26").Value = 5

# K3 (Days left), L3 (Total days = SUM(I3:I28)) and M3 (% done) are formulas
# that recalculate automatically from the I26 change above.

# Selection was moved to L26 before the workbook was saved.
$ws.Range("L26").Select()
